$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 7).Value = 42.05115733333333
$ws.Cells.Item(2, 8).Value = 126.153472
$ws.Cells.Item(2, 9).Value = 0.1594435451835853
$ws.Cells.Item(2, 10).Value = 0.1594435451835853
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.826632666666667
$ws.Cells.Item(2, 14).Value = 5.479898
$ws.Cells.Item(2, 15).Value = 0.4099391752648271
$ws.Cells.Item(2, 16).Value = 0.4099391752648271
$ws.Cells.Item(2, 17).Value = 76.81201765620622
$ws.Cells.Item(2, 18).Value = 691.308158905856
$ws.Cells.Item(2, 19).Value = 0.06536215541385915
$ws.Cells.Item(2, 20).Value = 0.06536215541385916

# Row 3
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 7).Value = 42.05115733333333
$ws.Cells.Item(3, 8).Value = 126.153472
$ws.Cells.Item(3, 9).Value = 0.1594435451835853
$ws.Cells.Item(3, 10).Value = 0.1594435451835853
$ws.Cells.Item(3, 13).Value = 2.304311333333333
$ws.Cells.Item(3, 14).Value = 6.912934
$ws.Cells.Item(3, 15).Value = 0.517141461870309
$ws.Cells.Item(3, 16).Value = 0.517141461870309
$ws.Cells.Item(3, 17).Value = 96.8989584229831
$ws.Cells.Item(3, 18).Value = 872.0906258068479
$ws.Cells.Item(3, 19).Value = 0.08245486804202395
$ws.Cells.Item(3, 20).Value = 0.08245486804202397

# Row 4
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 7).Value = 42.05115733333333
$ws.Cells.Item(4, 8).Value = 126.153472
$ws.Cells.Item(4, 9).Value = 0.1594435451835853
$ws.Cells.Item(4, 10).Value = 0.1594435451835853
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.3249186666666667
$ws.Cells.Item(4, 14).Value = 0.9747560000000001
$ws.Cells.Item(4, 15).Value = 0.07291936286486389
$ws.Cells.Item(4, 16).Value = 0.07291936286486389
$ws.Cells.Item(4, 17).Value = 13.66320597253689
$ws.Cells.Item(4, 18).Value = 122.968853752832
$ws.Cells.Item(4, 19).Value = 0.01162652172770218
$ws.Cells.Item(4, 20).Value = 0.01162652172770218

# Row 5
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 7).Value = 57.66057933333332
$ws.Cells.Item(5, 9).Value = 0.2186291119973147
$ws.Cells.Item(5, 10).Value = 0.2186291119973148
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.826632666666667
$ws.Cells.Item(5, 14).Value = 5.479898
$ws.Cells.Item(5, 15).Value = 0.4099391752648271
$ws.Cells.Item(5, 16).Value = 0.4099391752648271
$ws.Cells.Item(5, 17).Value = 105.3246977891915
$ws.Cells.Item(5, 18).Value = 947.9222801027239
$ws.Cells.Item(5, 19).Value = 0.08962463786106073
$ws.Cells.Item(5, 20).Value = 0.08962463786106074

# Row 6
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 7).Value = 57.66057933333332
$ws.Cells.Item(6, 9).Value = 0.2186291119973147
$ws.Cells.Item(6, 10).Value = 0.2186291119973148
$ws.Cells.Item(6, 13).Value = 2.304311333333333
$ws.Cells.Item(6, 14).Value = 6.912934
$ws.Cells.Item(6, 15).Value = 0.517141461870309
$ws.Cells.Item(6, 16).Value = 0.517141461870309
$ws.Cells.Item(6, 17).Value = 132.8679264443657
$ws.Cells.Item(6, 18).Value = 1195.811337999292
$ws.Cells.Item(6, 19).Value = 0.1130621785856988
$ws.Cells.Item(6, 20).Value = 0.1130621785856989

# Row 7
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 7).Value = 57.66057933333332
$ws.Cells.Item(7, 9).Value = 0.2186291119973147
$ws.Cells.Item(7, 10).Value = 0.2186291119973148
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.3249186666666667
$ws.Cells.Item(7, 14).Value = 0.9747560000000001
$ws.Cells.Item(7, 15).Value = 0.07291936286486389
$ws.Cells.Item(7, 16).Value = 0.07291936286486389
$ws.Cells.Item(7, 17).Value = 18.73499855621422
$ws.Cells.Item(7, 18).Value = 168.614987005928
$ws.Cells.Item(7, 19).Value = 0.01594229555055516
$ws.Cells.Item(7, 20).Value = 0.01594229555055516

# Row 8
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 7).Value = 99.15200299999999
$ws.Cells.Item(8, 8).Value = 297.456009
$ws.Cells.Item(8, 9).Value = 0.3759503393701321
$ws.Cells.Item(8, 10).Value = 0.3759503393701321
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.826632666666667
$ws.Cells.Item(8, 14).Value = 5.479898
$ws.Cells.Item(8, 15).Value = 0.4099391752648271
$ws.Cells.Item(8, 16).Value = 0.4099391752648271
$ws.Cells.Item(8, 17).Value = 181.1142876452313
$ws.Cells.Item(8, 18).Value = 1630.028588807082
$ws.Cells.Item(8, 19).Value = 0.1541167720619238
$ws.Cells.Item(8, 20).Value = 0.1541167720619238

# Row 9
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 7).Value = 99.15200299999999
$ws.Cells.Item(9, 8).Value = 297.456009
$ws.Cells.Item(9, 9).Value = 0.3759503393701321
$ws.Cells.Item(9, 10).Value = 0.3759503393701321
$ws.Cells.Item(9, 13).Value = 2.304311333333333
$ws.Cells.Item(9, 14).Value = 6.912934
$ws.Cells.Item(9, 15).Value = 0.517141461870309
$ws.Cells.Item(9, 16).Value = 0.517141461870309
$ws.Cells.Item(9, 17).Value = 228.4770842356006
$ws.Cells.Item(9, 18).Value = 2056.293758120406
$ws.Cells.Item(9, 19).Value = 0.1944195080925089
$ws.Cells.Item(9, 20).Value = 0.1944195080925089

# Row 10
$ws.Cells.Item(10, 4).Value = "Resolving-Mac"
$ws.Cells.Item(10, 7).Value = 99.15200299999999
$ws.Cells.Item(10, 8).Value = 297.456009
$ws.Cells.Item(10, 9).Value = 0.3759503393701321
$ws.Cells.Item(10, 10).Value = 0.3759503393701321
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.3249186666666667
$ws.Cells.Item(10, 14).Value = 0.9747560000000001
$ws.Cells.Item(10, 15).Value = 0.07291936286486389
$ws.Cells.Item(10, 16).Value = 0.07291936286486389
$ws.Cells.Item(10, 17).Value = 32.21633661208934
$ws.Cells.Item(10, 18).Value = 289.947029508804
$ws.Cells.Item(10, 19).Value = 0.02741405921569938
$ws.Cells.Item(10, 20).Value = 0.02741405921569939

# Row 11
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 7).Value = 64.87322933333333
$ws.Cells.Item(11, 8).Value = 194.619688
$ws.Cells.Item(11, 9).Value = 0.2459770034489679
$ws.Cells.Item(11, 10).Value = 0.2459770034489679
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.826632666666667
$ws.Cells.Item(11, 14).Value = 5.479898
$ws.Cells.Item(11, 15).Value = 0.4099391752648271
$ws.Cells.Item(11, 16).Value = 0.4099391752648271
$ws.Cells.Item(11, 17).Value = 118.4995598924249
$ws.Cells.Item(11, 18).Value = 1066.496039031824
$ws.Cells.Item(11, 19).Value = 0.1008356099279834
$ws.Cells.Item(11, 20).Value = 0.1008356099279834

# Row 12
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 7).Value = 64.87322933333333
$ws.Cells.Item(12, 8).Value = 194.619688
$ws.Cells.Item(12, 9).Value = 0.2459770034489679
$ws.Cells.Item(12, 10).Value = 0.2459770034489679
$ws.Cells.Item(12, 13).Value = 2.304311333333333
$ws.Cells.Item(12, 14).Value = 6.912934
$ws.Cells.Item(12, 15).Value = 0.517141461870309
$ws.Cells.Item(12, 16).Value = 0.517141461870309
$ws.Cells.Item(12, 17).Value = 149.4881175827324
$ws.Cells.Item(12, 18).Value = 1345.393058244592
$ws.Cells.Item(12, 19).Value = 0.1272049071500773
$ws.Cells.Item(12, 20).Value = 0.1272049071500773

# Row 13
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 7).Value = 64.87322933333333
$ws.Cells.Item(13, 8).Value = 194.619688
$ws.Cells.Item(13, 9).Value = 0.2459770034489679
$ws.Cells.Item(13, 10).Value = 0.2459770034489679
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.3249186666666667
$ws.Cells.Item(13, 14).Value = 0.9747560000000001
$ws.Cells.Item(13, 15).Value = 0.07291936286486389
$ws.Cells.Item(13, 16).Value = 0.07291936286486389
$ws.Cells.Item(13, 17).Value = 21.07852317734756
$ws.Cells.Item(13, 18).Value = 189.706708596128
$ws.Cells.Item(13, 19).Value = 0.01793648637090716
$ws.Cells.Item(13, 20).Value = 0.01793648637090717
